# Update imputed values in result_data_KNN sheet (Update Name of Algo)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4" = 8.075999999999999
    "B6" = 5.861
    "B7" = 5.395999999999999
    "E7" = 16.128
    "B8" = 5.959999999999999
    "E11" = 17.312
    "E12" = 17.638
    "E15" = 16.113
    "B16" = 5.621
    "B20" = 7.566
    "E20" = 16.291
    "B21" = 8.894
    "E21" = 16.821
    "E22" = 16.634
    "E23" = 16.527
    "B28" = 6.170999999999999
    "B29" = 5.306
    "E29" = 17.246
    "B30" = 6.112
    "B32" = 6.452999999999999
    "E34" = 16.77
    "B40" = 9.293000000000001
    "E42" = 16.539
    "E43" = 16.875
    "E44" = 16.376
    "E45" = 16.774
    "B46" = 6.382000000000001
    "E46" = 16.719
    "E50" = 16.548
    "B51" = 5.883999999999999
    "E51" = 16.881
    "B52" = 6.042
    "B57" = 5.188999999999999
    "E57" = 16.46
    "B59" = 5.125999999999999
    "B62" = 5.261
    "E65" = 17.109
    "B66" = 5.013
    "E66" = 17.384
    "E67" = 17.571
    "B73" = 6.813000000000001
    "B74" = 9.013000000000002
    "B77" = 5.751
    "E79" = 16.957
    "E84" = 16.606
    "E87" = 16.33
    "B92" = 5.198
    "E92" = 17.815
    "E97" = 16.794
    "B100" = 6.031000000000001
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
